$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Sistemata (!!!!) logica della calcolatrice" - new diary entry logging the
# implementation of the "potenza" (power) method, added as row 34.
$ws.Range("A34").Value = 43124                               # 24/01/2018
$ws.Range("B34").Value = "Mirko"
$ws.Range("C34").Value = "implementazione metodo potenza"
$ws.Range("D34").Value = 1/24                                # 1 hour

# Move the live selection to the next empty cell below the new row, matching
# where Excel leaves the cursor after the entry is typed in.
[void]$ws.Range("D35").Select()
